$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the description of dataset 3 with the new, longer text including
# special / HTML-like characters.
$ws.Range("G4").Value = "description of dataset 3, with speacial html l'ike > or & or < d'es fois"

# Update the active cell selection on the sheet.
$ws.Range("I7").Select()
